# Fruta / hortaliza, semanal
# Insert two new weekly report rows at the top of the data block (row 8),
# pushing the existing rows 8..49 down to rows 10..51 (dimension A1:T49 -> A1:T51),
# then populate the two newly inserted rows with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 8-49 down by two rows.
$ws.Rows("8:9").Insert()

# New row 8: Damasco, Modesto, Especial
$ws.Range("A8").Value = 8
$ws.Range("B8").Value = "Terminal La Palmera de La Serena"
$ws.Range("C8").Value = "Coquimbo"
$ws.Range("D8").Value = 44561
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = "Fruta"
$ws.Range("G8").Value = 100103
$ws.Range("H8").Value = "Frutos de hueso (carozo)"
$ws.Range("I8").Value = 100103003
$ws.Range("J8").Value = "Damasco"
$ws.Range("K8").Value = "Modesto"
$ws.Range("L8").Value = "Especial"
$ws.Range("M8").Value = 200
$ws.Range("N8").Value = 23000
$ws.Range("O8").Value = 24000
$ws.Range("P8").Value = 23500
$ws.Range("Q8").Value = "$/caja 18 kilos"
$ws.Range("R8").Value = "Región de O'Higgins"
$ws.Range("S8").Value = 1306
$ws.Range("T8").Value = 18

# New row 9: Damasco, Modesto, Primera
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Terminal La Palmera de La Serena"
$ws.Range("C9").Value = "Coquimbo"
$ws.Range("D9").Value = 44561
$ws.Range("E9").Value = 4
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100103
$ws.Range("H9").Value = "Frutos de hueso (carozo)"
$ws.Range("I9").Value = 100103003
$ws.Range("J9").Value = "Damasco"
$ws.Range("K9").Value = "Modesto"
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 200
$ws.Range("N9").Value = 20000
$ws.Range("O9").Value = 21000
$ws.Range("P9").Value = 20500
$ws.Range("Q9").Value = "$/caja 18 kilos"
$ws.Range("R9").Value = "Región de O'Higgins"
$ws.Range("S9").Value = 1139
$ws.Range("T9").Value = 18
